# Updated cryptos list on Fri Jul 12 20:24:18 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text so numeric-looking values such as
# "528.57" or "1.00" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.354.26"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.099.42"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "528.57"
$ws.Range("E5").Value = "  +0.84%  "

# Row 6 - Solana
$ws.Range("D6").Value = "137.20"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.102.44"

# Row 9 - XRP
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  +4.60%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").Value = "  +1.04%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.18%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.412"
$ws.Range("E12").Value = "  +4.35%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.36%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.649.40"
$ws.Range("E14").Value = "  +0.14%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "25.23"
$ws.Range("E15").Value = "  +0.13%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "57.517.80"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.111.17"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "5.98"
$ws.Range("E19").Value = "  +0.86%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "12.61"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "8.03"
$ws.Range("E21").Value = "  +2.08%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "360.11"
$ws.Range("E22").Value = "  +3.87%  "

# Row 23 - Dai
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "68.81"
$ws.Range("E24").Value = "  +0.91%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "0.502"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.72%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  +0.31%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0863"
$ws.Range("E28").Value = "  -4.55%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  -1.39%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").Value = "  -0.31%  "

# Row 31 - RenderToken
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  +0.35%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "21.30"
$ws.Range("E32").Value = "  +1.62%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "5.06"
$ws.Range("E33").Value = "  +3.05%  "

# Row 34 - Fetch.AI
$ws.Range("D34").Value = "1.13"
$ws.Range("E34").Value = "  -1.88%  "

# Row 35 - Monero
$ws.Range("D35").Value = "158.25"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  -1.72%  "

# Row 37 - EnergySwap
$ws.Range("D37").Value = "25.54"
$ws.Range("E37").Value = "  -1.41%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "1.26"
$ws.Range("E38").Value = "  +2.48%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "1.63"
$ws.Range("E39").Value = "  +2.12%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "0.0667"
$ws.Range("E40").Value = "  +0.65%  "

# Row 41 - Maker
$ws.Range("D41").Value = "2.480.30"
$ws.Range("E41").Value = "  +6.01%  "

# Row 42 - Filecoin/Mantle swap
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.696"
$ws.Range("E42").Value = "  -0.23%  "

# Row 43 - Mantle/Filecoin swap
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.98"
$ws.Range("E43").Value = "  -4.82%  "

# Row 44 - OKB
$ws.Range("D44").Value = "37.47"
$ws.Range("E44").Value = "  +3.02%  "

# Row 45 - FirstDigitalUSD
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46 - RenzoRestakedETH
$ws.Range("D46").Value = "3.148.66"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "0.0267"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "0.978"
$ws.Range("E48").Value = "  +3.11%  "

# Row 49 - Cosmos
$ws.Range("D49").Value = "6.03"
$ws.Range("E49").Value = "  +0.26%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "19.64"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "0.736"
$ws.Range("E51").Value = "  -3.21%  "

# Restore the default (Normal) style on the Price column so no stray
# number-format style is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
